$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B8: was a numeric 0, becomes the text "x" (text-formatted cell, new style w/ numFmtId 49)
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "x"

# B9: same change as B8
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "x"

# F14: value changes from 1 to 0
$ws.Range("F14").Value = 0

# Selection moves from B8 to F15
$ws.Range("F15").Select() | Out-Null
